$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure column D stays formatted as text so numeric-looking strings
# (e.g. "113.80", "3.00") are not coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "43.861.99"
$ws.Cells.Item(2, 5).Value = "  +0.12%  "

$ws.Cells.Item(3, 4).Value = "2.295.38"
$ws.Cells.Item(3, 5).Value = "  +0.26%  "

$ws.Cells.Item(4, 5).Value = "  +0.12%  "

$ws.Cells.Item(5, 4).Value = "113.80"
$ws.Cells.Item(5, 5).Value = "  +17.05%  "

$ws.Cells.Item(6, 4).Value = "269.74"
$ws.Cells.Item(6, 5).Value = "  -0.02%  "

$ws.Cells.Item(7, 5).Value = "  +0.38%  "

$ws.Cells.Item(8, 5).Value = "  +0.32%  "

$ws.Cells.Item(9, 5).Value = "  +1.84%  "

$ws.Cells.Item(10, 4).Value = "48.12"
$ws.Cells.Item(10, 5).Value = "  +5.96%  "

$ws.Cells.Item(11, 4).Value = "0.0944"
$ws.Cells.Item(11, 5).Value = "  +0.86%  "

$ws.Cells.Item(12, 4).Value = "9.06"
$ws.Cells.Item(12, 5).Value = "  +14.90%  "

$ws.Cells.Item(13, 5).Value = "  +0.29%  "

$ws.Cells.Item(14, 4).Value = "15.89"
$ws.Cells.Item(14, 5).Value = "  +0.25%  "

$ws.Cells.Item(15, 4).Value = "2.639.78"
$ws.Cells.Item(15, 5).Value = "  +0.20%  "

$ws.Cells.Item(16, 4).Value = "0.860"

$ws.Cells.Item(17, 4).Value = "2.295.13"
$ws.Cells.Item(17, 5).Value = "  +0.29%  "

$ws.Cells.Item(18, 4).Value = "43.751.10"
$ws.Cells.Item(18, 5).Value = "  -0.09%  "

$ws.Cells.Item(19, 5).Value = "  -1.08%  "

$ws.Cells.Item(20, 4).Value = "6.89"
$ws.Cells.Item(20, 5).Value = "  +11.07%  "

$ws.Cells.Item(21, 4).Value = "72.17"
$ws.Cells.Item(21, 5).Value = "  +0.09%  "

$ws.Cells.Item(22, 5).Value = "  -2.10%  "

$ws.Cells.Item(23, 4).Value = "3.00"
$ws.Cells.Item(23, 5).Value = "  +10.74%  "

$ws.Cells.Item(24, 4).Value = "232.85"
$ws.Cells.Item(24, 5).Value = "  -0.01%  "

$ws.Cells.Item(25, 4).Value = "9.68"
$ws.Cells.Item(25, 5).Value = "  +6.13%  "

$ws.Cells.Item(26, 4).Value = "0.999"
$ws.Cells.Item(26, 5).Value = "  -0.02%  "

$ws.Cells.Item(27, 4).Value = "11.64"
$ws.Cells.Item(27, 5).Value = "  +2.59%  "

$ws.Cells.Item(28, 5).Value = "  -1.17%  "

$ws.Cells.Item(29, 4).Value = "41.96"
$ws.Cells.Item(29, 5).Value = "  +8.46%  "

$ws.Cells.Item(30, 5).Value = "  -2.16%  "

$ws.Cells.Item(31, 5).Value = "  -0.80%  "

$ws.Cells.Item(32, 4).Value = "175.57"
$ws.Cells.Item(32, 5).Value = "  -0.07%  "

$ws.Cells.Item(33, 4).Value = "21.60"
$ws.Cells.Item(33, 5).Value = "  -0.93%  "

$ws.Cells.Item(34, 4).Value = "0.0924"
$ws.Cells.Item(34, 5).Value = "  +3.09%  "

$ws.Cells.Item(35, 4).Value = "5.72"
$ws.Cells.Item(35, 5).Value = "  +5.36%  "

$ws.Cells.Item(36, 5).Value = "  +0.01%  "

$ws.Cells.Item(37, 4).Value = "4.68"
$ws.Cells.Item(37, 5).Value = "  -0.57%  "

$ws.Cells.Item(38, 4).Value = "0.0364"
$ws.Cells.Item(38, 5).Value = "  +3.47%  "

$ws.Cells.Item(39, 5).Value = "  -0.02%  "

$ws.Cells.Item(40, 4).Value = "3.86"
$ws.Cells.Item(40, 5).Value = "  +12.02%  "

$ws.Cells.Item(41, 4).Value = "13.96"
$ws.Cells.Item(41, 5).Value = "  +13.85%  "

$ws.Cells.Item(42, 4).Value = "74.08"
$ws.Cells.Item(42, 5).Value = "  +14.91%  "

$ws.Cells.Item(43, 4).Value = "0.241"
$ws.Cells.Item(43, 5).Value = "  +1.93%  "

$ws.Cells.Item(44, 4).Value = "2.38"
$ws.Cells.Item(44, 5).Value = "  +2.86%  "

$ws.Cells.Item(45, 4).Value = "6.32"
$ws.Cells.Item(45, 5).Value = "  +21.43%  "

$ws.Cells.Item(46, 5).Value = "  +0.17%  "

$ws.Cells.Item(47, 5).Value = "  +2.92%  "

$ws.Cells.Item(48, 5).Value = "  +1.27%  "

$ws.Cells.Item(49, 4).Value = "0.0998"
$ws.Cells.Item(49, 5).Value = "  -2.64%  "

$ws.Cells.Item(50, 4).Value = "101.74"
$ws.Cells.Item(50, 5).Value = "  +3.31%  "

$ws.Cells.Item(51, 5).Value = "  +2.40%  "
